$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E23 was stored as text "524000"; the edit re-enters it as a real number.
$ws.Range("E23").Value() = 524000

# Append new row 24 (breakout of stock.yaml data).
$ws.Range("A24").Value() = "19/06/2024 09:47:02"
$ws.Range("B24").Value() = 1
$ws.Range("C24").Value() = "POONAWALLA"
$ws.Range("D24").Value() = "Poonawalla Fincorp Ltd"

# E24 keeps the bsecode as text (leading apostrophe forces text entry),
# then strip the resulting quote-prefix style so no extra formatting sticks.
$ws.Range("E24").Value() = "'524000"
$ws.Range("E24").Style = "Normal"

$ws.Range("F24").Value() = -2.2
$ws.Range("G24").Value() = 424
$ws.Range("H24").Value() = 6147053
